$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("D10").Value = 56
$ws.Range("E10").Value = 56
$ws.Range("C11").Value = "'-1"
$ws.Range("C11").Style = "incorrectStyle"
$ws.Range("G15").Value = "Student Ans"
$ws.Range("G15").Style = "mtitleStyle"
$ws.Range("H15").Value = "Correct Ans"
$ws.Range("H15").Style = "mtitleStyle"
$ws.Range("G16").Value = ""
$ws.Range("G16").Style = "normalStyle"
$ws.Range("H16").Value = "Option A"
$ws.Range("H16").Style = "absoluteStyle"
$ws.Range("G17").Value = ""
$ws.Range("G17").Style = "normalStyle"
$ws.Range("H17").Value = "Option D"
$ws.Range("H17").Style = "absoluteStyle"
$ws.Range("G18").Value = ""
$ws.Range("G18").Style = "normalStyle"
$ws.Range("H18").Value = "Option D"
$ws.Range("H18").Style = "absoluteStyle"
$ws.Range("D19").Value = ""
$ws.Range("D19").Style = "normalStyle"
$ws.Range("E19").Value = "Option A"
$ws.Range("E19").Style = "absoluteStyle"
$ws.Range("G19").Value = ""
$ws.Range("G19").Style = "normalStyle"
$ws.Range("H19").Value = "Option A"
$ws.Range("H19").Style = "absoluteStyle"
$ws.Range("D20").Value = ""
$ws.Range("D20").Style = "normalStyle"
$ws.Range("E20").Value = "Option D"
$ws.Range("E20").Style = "absoluteStyle"
$ws.Range("G20").Value = ""
$ws.Range("G20").Style = "normalStyle"
$ws.Range("H20").Value = "Option C"
$ws.Range("H20").Style = "absoluteStyle"
$ws.Range("D21").Value = ""
$ws.Range("D21").Style = "normalStyle"
$ws.Range("E21").Value = "Option B"
$ws.Range("E21").Style = "absoluteStyle"
$ws.Range("G21").Value = ""
$ws.Range("G21").Style = "normalStyle"
$ws.Range("H21").Value = "Option D"
$ws.Range("H21").Style = "absoluteStyle"
$ws.Range("D22").Value = ""
$ws.Range("D22").Style = "normalStyle"
$ws.Range("E22").Value = "Option C"
$ws.Range("E22").Style = "absoluteStyle"
$ws.Range("D23").Value = ""
$ws.Range("D23").Style = "normalStyle"
$ws.Range("E23").Value = "Option B"
$ws.Range("E23").Style = "absoluteStyle"
$ws.Range("D24").Value = ""
$ws.Range("D24").Style = "normalStyle"
$ws.Range("E24").Value = "Option C"
$ws.Range("E24").Style = "absoluteStyle"
$ws.Range("D25").Value = ""
$ws.Range("D25").Style = "normalStyle"
$ws.Range("E25").Value = "Option D"
$ws.Range("E25").Style = "absoluteStyle"
$ws.Range("D26").Value = ""
$ws.Range("D26").Style = "normalStyle"
$ws.Range("E26").Value = "Option D"
$ws.Range("E26").Style = "absoluteStyle"
$ws.Range("D27").Value = ""
$ws.Range("D27").Style = "normalStyle"
$ws.Range("E27").Value = "Option A"
$ws.Range("E27").Style = "absoluteStyle"
$ws.Range("D28").Value = ""
$ws.Range("D28").Style = "normalStyle"
$ws.Range("E28").Value = "Option A"
$ws.Range("E28").Style = "absoluteStyle"
$ws.Range("D29").Value = ""
$ws.Range("D29").Style = "normalStyle"
$ws.Range("E29").Value = "Option C"
$ws.Range("E29").Style = "absoluteStyle"
$ws.Range("D30").Value = ""
$ws.Range("D30").Style = "normalStyle"
$ws.Range("E30").Value = "Option A"
$ws.Range("E30").Style = "absoluteStyle"
$ws.Range("D31").Value = ""
$ws.Range("D31").Style = "normalStyle"
$ws.Range("E31").Value = "Option D"
$ws.Range("E31").Style = "absoluteStyle"
$ws.Range("D32").Value = ""
$ws.Range("D32").Style = "normalStyle"
$ws.Range("E32").Value = "Option D"
$ws.Range("E32").Style = "absoluteStyle"
$ws.Range("D33").Value = ""
$ws.Range("D33").Style = "normalStyle"
$ws.Range("E33").Value = "Option B"
$ws.Range("E33").Style = "absoluteStyle"
$ws.Range("D34").Value = ""
$ws.Range("D34").Style = "normalStyle"
$ws.Range("E34").Value = "Option D"
$ws.Range("E34").Style = "absoluteStyle"
$ws.Range("D35").Value = ""
$ws.Range("D35").Style = "normalStyle"
$ws.Range("E35").Value = "Option C"
$ws.Range("E35").Style = "absoluteStyle"
$ws.Range("D36").Value = ""
$ws.Range("D36").Style = "normalStyle"
$ws.Range("E36").Value = "Option D"
$ws.Range("E36").Style = "absoluteStyle"
$ws.Range("D37").Value = ""
$ws.Range("D37").Style = "normalStyle"
$ws.Range("E37").Value = "Option B"
$ws.Range("E37").Style = "absoluteStyle"
$ws.Range("D38").Value = ""
$ws.Range("D38").Style = "normalStyle"
$ws.Range("E38").Value = "Option D"
$ws.Range("E38").Style = "absoluteStyle"
$ws.Range("D39").Value = ""
$ws.Range("D39").Style = "normalStyle"
$ws.Range("E39").Value = "Option A"
$ws.Range("E39").Style = "absoluteStyle"
$ws.Range("D40").Value = ""
$ws.Range("D40").Style = "normalStyle"
$ws.Range("E40").Value = "Option A"
$ws.Range("E40").Style = "absoluteStyle"
